$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 6-12 hold the "Energie"/"Industriemetalle" commodity rows.
# The "Diesel" row is removed (rows below shift up by one), the
# "Energie" category is renamed to "Energie u. weiteres", and a new
# "Kakao" row is appended at the bottom with the same indicator value
# that "Zink" used to carry before the shift.

$ws.Range("A6").Value = "Benzin"
$ws.Range("B6").Value = "Energie u. weiteres"
$ws.Range("E6").Value = 3064

$ws.Range("A7").Value = "Erdgas"
$ws.Range("B7").Value = "Energie u. weiteres"
$ws.Range("E7").Value = 2773

$ws.Range("A8").Value = "Rohöl"
$ws.Range("B8").Value = "Energie u. weiteres"
$ws.Range("E8").Value = 3480

$ws.Range("A9").Value = "Aluminium"
$ws.Range("B9").Value = "Industriemetalle"
$ws.Range("E9").Value = 3054

$ws.Range("A10").Value = "Kupfer"
$ws.Range("B10").Value = "Industriemetalle"
$ws.Range("E10").Value = 2943

$ws.Range("A11").Value = "Zink"
$ws.Range("B11").Value = "Industriemetalle"
$ws.Range("E11").Value = 3053

$ws.Range("A12").Value = "Kakao"
$ws.Range("B12").Value = "Energie u. weiteres"
$ws.Range("D12").Value = "Rohstoff Indikator"
$ws.Range("E12").Value = 3053

# Update the selected cell shown in the saved view.
$ws.Range("I16").Select()
